# The "exceed" counters in columns B/C/D (L1/L2/L3) used to be plain
# numbers. This change switches the cells that represent an exceeded
# count (e.g. "6+", "17+") - and, for consistency, their siblings in the
# same rows - to store their value as text instead of a number.
#
# Writing a numeric-looking string (e.g. "0", "1") straight into
# Range.Value gets auto-coerced back into a number by Excel, so those
# cells are written as a text formula (leading apostrophe) which forces
# text storage, and the style is immediately reset back to "Normal" so
# no formatting actually changes on the cell - only the underlying
# value type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# Row 2 (L1=6+)
$ws.Range("B2").Value = "6+"

# Row 3 (L1=0, L2=0+)
Set-TextValue $ws.Range("B3") "0"
Set-TextValue $ws.Range("C3") "0+"

# Row 4 (L1=0, L2=0+)
Set-TextValue $ws.Range("B4") "0"
Set-TextValue $ws.Range("C4") "0+"

# Row 5 (L1=0, L2=1, L3=17+)
Set-TextValue $ws.Range("B5") "0"
Set-TextValue $ws.Range("C5") "1"
$ws.Range("D5").Value = "17+"
